$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.773.15"
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").Value = "2.278.47"
$ws.Range("E3").Value = "  +1.29%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.00"
$ws.Range("E5").Value = "  +0.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.634"
$ws.Range("E6").Value = "  +0.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.32"
$ws.Range("E7").Value = "  +7.13%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.644"
$ws.Range("E9").Value = "  -2.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.41"
$ws.Range("E10").Value = "  +1.23%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0977"
$ws.Range("E11").Value = "  +1.53%  "

$ws.Range("E12").Value = "  -0.41%  "

$ws.Range("E13").Value = "  +2.34%  "

$ws.Range("D14").Value = "2.618.23"
$ws.Range("E14").Value = "  +1.46%  "

$ws.Range("E15").Value = "  +2.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.869"
$ws.Range("E16").Value = "  -1.06%  "

$ws.Range("D17").Value = "2.287.81"
$ws.Range("E17").Value = "  +0.43%  "

$ws.Range("D18").Value = "42.631.47"
$ws.Range("E18").Value = "  -0.14%  "

$ws.Range("D19").Value = "0.0₃0996"
$ws.Range("E19").Value = "  +1.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.20"
$ws.Range("E20").Value = "  -1.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.40"
$ws.Range("E21").Value = "  -0.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "236.28"
$ws.Range("E22").Value = "  +0.62%  "

$ws.Range("E23").Value = "  +6.10%  "

$ws.Range("E24").Value = "  -1.55%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.26"
$ws.Range("E26").Value = "  -1.34%  "

$ws.Range("E27").Value = "  -1.29%  "

$ws.Range("E28").Value = "  -2.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.55"
$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("E30").Value = "  +1.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0873"
$ws.Range("E31").Value = "  +9.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.33"
$ws.Range("E32").Value = "  -1.78%  "

$ws.Range("E33").Value = "  +0.34%  "

$ws.Range("E34").Value = "  +2.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.128"
$ws.Range("E35").Value = "  +2.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.54"
$ws.Range("E36").Value = "  +2.79%  "

$ws.Range("E37").Value = "  +1.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0305"
$ws.Range("E38").Value = "  -5.40%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.49"
$ws.Range("E39").Value = "  +9.16%  "

$ws.Range("E40").Value = "  -0.11%  "

$ws.Range("E41").Value = "  +1.53%  "

$ws.Range("E42").Value = "  +3.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "61.50"
$ws.Range("E43").Value = "  -0.99%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.93"
$ws.Range("E44").Value = "  -0.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "105.55"
$ws.Range("E45").Value = "  +11.76%  "

$ws.Range("E46").Value = "  -3.02%  "

$ws.Range("E47").Value = "  -0.69%  "

$ws.Range("E48").Value = "  +0.25%  "

$ws.Range("E49").Value = "  +0.07%  "

$ws.Range("E50").Value = "  -1.42%  "

$ws.Range("E51").Value = "  -0.58%  "
